$d = $word.ActiveDocument

function New-PkgXml($innerBodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Paragraph 1: "users(...)" -> "users (...)" with proofErr stripped and runs re-split ---
$para1Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Sansinterligne"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>u</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sers</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>email</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, last_name, first_name, password, sold, role)</w:t></w:r>' +
  '</w:p></w:body>'
$xml1 = New-PkgXml $para1Body
$d.Paragraphs.Item(1).Range.InsertXML($xml1)

# --- Paragraph 2: "transactions(...)" -> "transactions (...)" with new columns, proofErr stripped ---
$para2Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Sansinterligne"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>transactions (</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>id_transaction</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, description, amount</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, #email</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, #email_friend</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
  '</w:p></w:body>'
$xml2 = New-PkgXml $para2Body
$d.Paragraphs.Item(2).Range.InsertXML($xml2)

# --- Paragraphs 3+4: "add_friend(...)" reworked, and "user_transaction(...)" paragraph removed entirely ---
$para3Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Sansinterligne"/><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>add_friend</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="words"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">(#email </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>,</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="words"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="words"/><w:lang w:val="en-US"/></w:rPr><w:t>#email_friend</w:t></w:r>' +
  '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
  '</w:p></w:body>'
$xml3 = New-PkgXml $para3Body

$p3Range = $d.Paragraphs.Item(3).Range
$p4Range = $d.Paragraphs.Item(4).Range
$combined = $d.Range($p3Range.Start, $p4Range.End)
$combined.InsertXML($xml3)
